# plan.xlsx update: refresh the product-code list on Planilha1
# - remove the stale codes (C-2062, C-2044, C-3308, C-3334)
# - keep the remaining codes (C-3377, C-1125, C-1478)
# - append the newly-tracked codes at the bottom of the list

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four stale rows (rows 2-5: C-2062, C-2044, C-3308, C-3334).
# The remaining rows shift up automatically.
$ws.Range("A2:A5").EntireRow.Delete() | Out-Null

# Append the new product codes below the existing data.
$newCodes = @("C-1528", "C-1526", "C-1500", "L-663", "L-733", "L-533X", "L-533", "L-562")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($i = 0; $i -lt $newCodes.Count; $i++) {
    $ws.Cells.Item($lastRow + 1 + $i, 1).Value = $newCodes[$i]
}

# Keep the hidden _FilterDatabase defined name in sync with the new range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Planilha1!_FilterDatabase") {
        $n.RefersTo = "=Planilha1!`$A`$1:`$A`$109"
    }
}
